$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30/31: swap PancakeSwap/Bittensor data with updated price/volume values
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "547.29"
$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").Value = "  +8.88%  "

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.335.28"
$ws.Range("E2").Value = "  +3.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.70"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.98"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.91"
$ws.Range("E6").Value = "  +5.39%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +11.34%  "
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.38"
$ws.Range("E13").Value = "  +6.86%  "
$ws.Range("E14").Value = "  +22.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.121.32"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.153.53"
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.662.15"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.70"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("E19").Value = "  +5.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.25"
$ws.Range("E20").Value = "  +4.89%  "
$ws.Range("E21").Value = "  +8.49%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.28"
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.45"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0970"
$ws.Range("E29").Value = "  +14.66%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("E34").Value = "  +6.23%  "
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.60"
$ws.Range("E37").Value = "  +6.60%  "
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.53"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.83"
$ws.Range("E42").Value = "  +8.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.53"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0619"
$ws.Range("E45").Value = "  +7.43%  "
$ws.Range("E46").Value = "  +9.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.28"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.656"
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("E49").Value = "  +7.16%  "
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("E51").Value = "  +4.15%  "
